# Update the cached "datetimeFigureOut" footer field text from 2020-06-28
# to 2020-06-29 on every Date Placeholder shape across the slide master
# and all slide (custom) layouts.

$p = $ppt.ActivePresentation

$oldDate = "2020-06-28"
$newDate = "2020-06-29"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# All custom (slide) layouts under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $lay = $layouts.Item($j)
    Update-DatePlaceholder $lay.Shapes
}

Write-Host "Date placeholders updated."
